# Update "想去人数" (want-to-go count) values that changed between scrapes.
# Sheet "展览" and Sheet "全部类型" each mirror the same rows:
#   F2: 332 -> 333
#   F5: 285 -> 287

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 333
    $ws.Range("F5").Value = 287
}
